$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gender & Race")

# Insert a new row 6 (pushes the existing rows 6-51 down to 7-52), mirroring
# the existing "date / link text" rows above it (rows 1-5).
$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = 43513
$ws.Range("B6").Value = "http://share.streamlit.io/0.25.0-cdyb/index.html?id=2vHQ1bySr6a1FnzM6x14De"

$ws.Hyperlinks.Add($ws.Range("B6"), "http://share.streamlit.io/0.25.0-cdyb/index.html?id=2vHQ1bySr6a1FnzM6x14De")

# Copy the style (date + link-text) from the row above so row 6 matches rows 1-5.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("B5").Copy()
$ws.Range("B6").PasteSpecial(-4122)

$ws.Range("D3").Select()

$ws.Tab.Activate()
